$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Revert "update meeting log": clear the row 23 entries that were added
# (B23:F23), leaving only the row's "#" index (A23) and trailing G23 cell.
$ws.Range("B23:F23").ClearContents()

# Restore the prior selection on the sheet (was moved to B23 by the
# original edit; revert moves it to F22).
$ws.Range("F22").Select()
